# Fonds de solidarite - add 2020-08-06 data
# Updates nombre_aides (col C) and montant_total (col D) for the rows whose
# classe_effectif counts/amounts changed with the new data pull.
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the sheet's existing inline-string / text-typed cell layout,
# e.g. "235826.00" instead of being auto-coerced to the number 235826).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 33; C = "89";  D = "235826.00" },
    @{ Row = 34; C = "491"; D = "1459400.41" },
    @{ Row = 35; C = "198"; D = "937647.11" },
    @{ Row = 36; C = "71";  D = "383474.00" },
    @{ Row = 39; C = "30";  D = "74330.00" },
    @{ Row = 45; C = "22";  D = "92621.84" },
    @{ Row = 46; C = "59";  D = "263191.74" },
    @{ Row = 47; C = "36";  D = "204937.00" },
    @{ Row = 48; C = "23";  D = "153697.00" },
    @{ Row = 72; C = "8";   D = "39000.00" },
    @{ Row = 79; C = "217"; D = "551826.09" },
    @{ Row = 80; C = "838"; D = "2555291.11" },
    @{ Row = 81; C = "312"; D = "1198440.79" },
    @{ Row = 82; C = "103"; D = "462484.52" },
    @{ Row = 92; C = "388"; D = "1110260.67" },
    @{ Row = 94; C = "47";  D = "221347.01" }
)

foreach ($u in $updates) {
    $cC = $ws.Cells.Item($u.Row, 3)
    $cC.Value = "'" + $u.C
    $cC.Style = "Normal"

    $cD = $ws.Cells.Item($u.Row, 4)
    $cD.Value = "'" + $u.D
    $cD.Style = "Normal"
}
